$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 0.2
$ws.Range("F2").Value = 1234567

$ws.Range("E3").Value = 0.3
$ws.Range("F3").Value = 123456

$ws.Range("E4").Value = 0.4
$ws.Range("F4").Value = 12345

$ws.Range("E5").Value = 0.1
$ws.Range("F5").Value = 1234

$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 123

$ws.Range("E7").Value = 0.4
$ws.Range("F7").Value = 12345

$ws.Range("E8").Value = 0.6
$ws.Range("F8").Value = 123

$ws.Range("E9").Value = 0.7
$ws.Range("F9").Value = 123678

$ws.Range("E10").Value = 0.3
$ws.Range("F10").Value = 123

$ws.Range("E11").Value = 0.66
$ws.Range("F11").Value = 56789

$ws.Range("E12").Value = 0.34
$ws.Range("F12").Value = 123

$ws.Range("E13").Value = 0.4
$ws.Range("F13").Value = 469

$ws.Range("E14").Value = 0.6
$ws.Range("F14").Value = 136

$ws.Range("F15").Value = 56789

$ws.Range("E16").Value = 0.59
$ws.Range("F16").Value = 987
$ws.Range("G16").Value = $True

$ws.Range("E17").Value = 0.4
$ws.Range("F17").Value = 345

$ws.Range("E18").Value = 0.2
$ws.Range("F18").Value = 567

$ws.Range("E19").Value = 0.4
$ws.Range("F19").Value = 765

$ws.Range("E20").Value = 0.8

$ws.Range("E21").Value = 0.2
$ws.Range("F21").Value = 434433

$ws.Range("E22").Value = 0.2
$ws.Range("F22").Value = 434343

$ws.Range("E23").Value = 0.4
$ws.Range("F23").Value = 434343

$ws.Range("E24").Value = 0.4
$ws.Range("F24").Value = 434

$ws.Range("E25").Value = 0.56
$ws.Range("F25").Value = 7658

$ws.Range("E26").Value = 0.44
$ws.Range("F26").Value = 874678

$ws.Range("E27").Value = 0.4
$ws.Range("F27").Value = 232323
$ws.Range("G27").Value = $True

$ws.Range("E28").Value = 0.2
$ws.Range("F28").Value = 4564

$ws.Range("E29").Value = 0.2
$ws.Range("F29").Value = 5456

$ws.Range("E30").Value = 0.2
$ws.Range("F30").Value = 5455

$ws.Range("E31").Value = 0.2
$ws.Range("F31").Value = 45454

$ws.Range("E32").Value = 0.8
$ws.Range("F32").Value = 45454

$ws.Range("E33").Value = 0.9
$ws.Range("F33").Value = 6666

$ws.Range("E34").Value = 0.1
$ws.Range("F34").Value = 443

$ws.Range("B15").Select()
